$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (N1:Q1) - bold/centered/bordered like existing headers ---
$ws.Range("A1").Copy() | Out-Null
$ws.Range("N1:Q1").PasteSpecial(-4122) | Out-Null

$ws.Range("N1").Value = 'bot_generated_search_string'
$ws.Range("O1").Value = 'ss_words'
$ws.Range("P1").Value = 'bot_gen_ss_words'
$ws.Range("Q1").Value = 'jaccard_similarity'

# --- Data rows 2-26: N (bot_generated_search_string), O (ss_words), P (bot_gen_ss_words), Q (jaccard_similarity) ---
$ws.Range("N2").Value = '("software maintenance" OR "software evolution" OR "legacy systems") AND ("globally distributed" OR "distributed teams" OR "offshore development" OR "global software engineering") AND ("challenges" OR "problems" OR "issues" OR "difficulties" OR "obstacles") AND ("mitigation" OR "resolution" OR "strategies" OR "approaches" OR "solutions" OR "best practices")'
$ws.Range("O2").Value = '[''change'', ''distributed software development'', ''distributed software engineering'', ''distributed software project'', ''global enterprise resource planning erp software'', ''global enterprise resource planning software'', ''global software development'', ''global software engineering'', ''global software project'', ''maintain'', ''maintenance'', ''software offshore'', ''software outsource'', ''support'', ''upgrade'']'
$ws.Range("P2").Value = '[''approaches'', ''best practices'', ''challenges'', ''difficulties'', ''distributed teams'', ''global software engineering'', ''globally distributed'', ''issues'', ''legacy systems'', ''mitigation'', ''obstacles'', ''offshore development'', ''problems'', ''resolution'', ''software evolution'', ''software maintenance'', ''solutions'', ''strategies'']'
$ws.Range("Q2").Value = 0.03125
$ws.Range("N3").Value = '("success factors" OR "critical success factors" OR "enablers" OR "drivers") AND ("chief data officer" OR "CDO" OR "data leader" OR "data governance") AND ("literature review" OR "expert opinion" OR "questionnaire" OR "survey") AND ("prioritization" OR "taxonomy" OR "classification" OR "framework")'
$ws.Range("O3").Value = '[''application service'', ''asp"and"cloud computing'', ''aspects'', ''cloud offering'', ''cloud platform'', ''cloud provider'', ''cloud service'', ''collaborative software development"and"iaas'', ''drivers'', ''elements'', ''factors'', ''geographically distributed development'', ''global software development'', ''infrastructure as a service'', ''it service'', ''items'', ''motivators'', ''multisite development'', ''offshore development'', ''paas'', ''platform as a service'', ''saas'', ''software as a service'', ''success factors'', ''variables"and"outsourcing'', ''xaas'']'
$ws.Range("P3").Value = '[''cdo'', ''chief data officer'', ''classification'', ''critical success factors'', ''data governance'', ''data leader'', ''drivers'', ''enablers'', ''expert opinion'', ''framework'', ''literature review'', ''prioritization'', ''questionnaire'', ''success factors'', ''survey'', ''taxonomy'']'
$ws.Range("Q3").Value = 0.05
$ws.Range("N4").Value = '("model transformation" OR "model transformations") AND ("design pattern" OR "design patterns" OR "pattern language") AND ("MT development" OR "model transformation development") AND (practice OR application OR usage) AND (categories OR types OR classification) AND (explicit OR recognized) AND (benefits OR advantages OR outcomes) AND (adoption OR trend OR evolution) AND (languages OR frameworks OR tools)'
$ws.Range("O4").Value = '[''model transformation'']'
$ws.Range("P4").Value = '[''adoption'', ''advantages'', ''application'', ''benefits'', ''categories'', ''classification'', ''design pattern'', ''design patterns'', ''evolution'', ''explicit'', ''frameworks'', ''languages'', ''model transformation'', ''model transformation development'', ''model transformations'', ''mt development'', ''outcomes'', ''pattern language'', ''practice'', ''recognized'', ''tools'', ''trend'', ''types'', ''usage'']'
$ws.Range("Q4").Value = 0.04166666666666666
$ws.Range("N5").Value = '("blockchain governance" OR "decentralized governance" OR "on-chain governance" OR "distributed governance") AND ("building blocks" OR "components" OR "elements" OR "framework" OR "architecture")'
$ws.Range("O5").Value = '[''autonomous organization"and governance'', ''blockchain ordistributed'', ''decentralizedand ledger'', ''ecosystem'', ''management'', ''platform'']'
$ws.Range("P5").Value = '[''architecture'', ''blockchain governance'', ''building blocks'', ''components'', ''decentralized governance'', ''distributed governance'', ''elements'', ''framework'', ''on-chain governance'']'
$ws.Range("Q5").Value = 0
$ws.Range("N6").Value = '("mortality compression" OR "MCR") AND ("approaches" OR "methods" OR "strategies") AND ("evaluation" OR "assessment" OR "outcomes") AND ("conclusions" OR "findings" OR "results")'
$ws.Range("O6").Value = '[''code inspection'', ''code review'', ''formal inspection'', ''software inspection'']'
$ws.Range("P6").Value = '[''approaches'', ''assessment'', ''conclusions'', ''evaluation'', ''findings'', ''mcr'', ''methods'', ''mortality compression'', ''outcomes'', ''results'', ''strategies'']'
$ws.Range("Q6").Value = 0
$ws.Range("N7").Value = ""
$ws.Range("O7").Value = '[''“cognitive effectiveness”'', ''“diagram”'', ''“language”'', ''“modeling”'', ''“notation”'', ''“physics of notations”'', ''“visual”'']'
$ws.Range("P7").Value = '[''nan'']'
$ws.Range("Q7").Value = 0
$ws.Range("N8").Value = '(("semi-automatic configuration" OR "automated configuration" OR "assisted configuration") AND ("product line" OR "product family" OR "software product line") AND ("techniques" OR "methods" OR "approaches") AND ("evaluation" OR "validation" OR "assessment") AND ("challenges" OR "limitations" OR "open issues"))'
$ws.Range("O8").Value = '[''attribute'', ''confugyration'', ''feature selection'', ''non-functional'', ''preference'', ''product family'', ''product line'', ''product selection'', ''quality'', ''requirement'', ''system family'']'
$ws.Range("P8").Value = '[''approaches'', ''assessment'', ''assisted configuration'', ''automated configuration'', ''challenges'', ''evaluation'', ''limitations'', ''methods'', ''open issues'', ''product family'', ''product line'', ''semi-automatic configuration'', ''software product line'', ''techniques'', ''validation'']'
$ws.Range("Q8").Value = 0.08333333333333333
$ws.Range("N9").Value = '("Deep Learning" OR "Deep Neural Networks" OR "Machine Learning") AND ("Software Engineering" OR "Software Development" OR "Software Maintenance") AND (Tasks OR Applications) AND (Data Extraction OR Data Preprocessing) AND (Model Architecture OR Learning Algorithms) AND (Performance Evaluation OR Benchmarks) AND (Reproducibility OR Replication)'
$ws.Range("O9").Value = '[''“deep”'', ''“learning”'', ''“neural”'']'
$ws.Range("P9").Value = '[''applications'', ''benchmarks'', ''data extraction'', ''data preprocessing'', ''deep learning'', ''deep neural networks'', ''learning algorithms'', ''machine learning'', ''model architecture'', ''performance evaluation'', ''replication'', ''reproducibility'', ''software development'', ''software engineering'', ''software maintenance'', ''tasks'']'
$ws.Range("Q9").Value = 0
$ws.Range("N10").Value = '("software trust" OR "SECO trust" OR "security trust") AND (definition OR concept OR meaning) AND (("software product" OR "software version" OR "software package manager" OR "software organization" OR "software engineer") AND (selection OR choosing OR decision-making) AND (trust OR trustworthiness OR reliability OR security OR safety OR integrity OR reputation OR assurance))'
$ws.Range("O10").Value = '[''compoonent'', ''credibility'', ''dependency'', ''developer'', ''management'', ''package'', ''procenance'', ''reputation'', ''software'', ''software ecosystem'', ''trust'', ''uncertainty'']'
$ws.Range("P10").Value = '[''assurance'', ''choosing'', ''concept'', ''decision-making'', ''definition'', ''integrity'', ''meaning'', ''reliability'', ''reputation'', ''safety'', ''seco trust'', ''security'', ''security trust'', ''selection'', ''software engineer'', ''software organization'', ''software package manager'', ''software product'', ''software trust'', ''software version'', ''trust'', ''trustworthiness'']'
$ws.Range("Q10").Value = 0.0625
$ws.Range("N11").Value = '("data mining" OR "machine learning" OR "text mining" OR "sentiment analysis" OR "topic modeling" OR "classification" OR "clustering" OR "regression") AND ("app store" OR "Google Play" OR "Apple App Store") AND ("software review" OR "app review") AND ("domain adaptation" OR "transfer learning" OR "cross-domain" OR "domain shift" OR "contextual variation") AND ("review helpfulness" OR "review quality" OR "review credibility" OR "review informativeness" OR "accuracy" OR "relevance" OR "completeness" OR "timeliness" OR "specificity") AND ("spam detection" OR "fraudulent review" OR "fake review" OR "deceptive review") AND ("feature extraction" OR "aspect extraction" OR "information extraction" OR "entity recognition") AND ("software feature" OR "app feature" OR "performance" OR "usability" OR "security" OR "functionality")'
$ws.Range("O11").Value = '[''analysis'', ''app'', ''application'', ''apps'', ''appstore'', ''bug'', ''comment'', ''complain'', ''data'', ''ecosystem'', ''expectation'', ''feature'', ''helpfulness'', ''issue'', ''market'', ''mining'', ''mobile'', ''online'', ''opinion'', ''processing'', ''quality'', ''rating'', ''request'', ''requirement'', ''review'', ''sentiment'', ''software'', ''store'', ''text'', ''usefulness'', ''user'', ''vocabulary'']'
$ws.Range("P11").Value = '[''accuracy'', ''app feature'', ''app review'', ''app store'', ''apple app store'', ''aspect extraction'', ''classification'', ''clustering'', ''completeness'', ''contextual variation'', ''cross-domain'', ''data mining'', ''deceptive review'', ''domain adaptation'', ''domain shift'', ''entity recognition'', ''fake review'', ''feature extraction'', ''fraudulent review'', ''functionality'', ''google play'', ''information extraction'', ''machine learning'', ''performance'', ''regression'', ''relevance'', ''review credibility'', ''review helpfulness'', ''review informativeness'', ''review quality'', ''security'', ''sentiment analysis'', ''software feature'', ''software review'', ''spam detection'', ''specificity'', ''text mining'', ''timeliness'', ''topic modeling'', ''transfer learning'', ''usability'']'
$ws.Range("Q11").Value = 0
$ws.Range("N12").Value = ""
$ws.Range("O12").Value = '[''architecting'', ''architectural'', ''architecture'', ''diagram'', ''graphic'', ''graphical'', ''picture'', ''structure'', ''visual'', ''visualization'', ''visualize'', ''visualizing'']'
$ws.Range("P12").Value = '[''nan'']'
$ws.Range("Q12").Value = 0
$ws.Range("N13").Value = '("UML" OR "Unified Modeling Language") AND ("consistency" OR "validation" OR "verification") AND ("diagram" OR "model") AND ("formal method" OR "formal technique" OR "non-formal method" OR "non-formal technique") AND ("version" OR "revision") AND ("check" OR "rule" OR "constraint")'
$ws.Range("O13").Value = '[''“consistency”'', ''“consistency”or “model”'', ''“inconsistency”'', ''“inconsistency”'', ''“management”'', ''“model”'']'
$ws.Range("P13").Value = '[''check'', ''consistency'', ''constraint'', ''diagram'', ''formal method'', ''formal technique'', ''model'', ''non-formal method'', ''non-formal technique'', ''revision'', ''rule'', ''uml'', ''unified modeling language'', ''validation'', ''verification'', ''version'']'
$ws.Range("Q13").Value = 0
$ws.Range("N14").Value = '("Search Based Software Testing" OR SBST) AND ("Mutation Testing" OR "Fault Injection") AND (Meta-heuristics OR "Genetic Algorithm" OR "Simulated Annealing" OR "Evolutionary Algorithm") AND ("Fitness Functions") AND ("Programming Languages")'
$ws.Range("O14").Value = '[''ant colony'', ''evelutionary'', ''genetic algorithms'', ''genetic programming'', ''heuristic'', ''hill-climbing'', ''meta-heuristic'', ''metaheuristic'', ''mutation analysis'', ''mutation based test'', ''mutation test'', ''mutation testing'', ''mutation-based test'', ''oprimization'', ''program mutation'', ''search based'', ''search-based'', ''simulated annealing'', ''tabu search'']'
$ws.Range("P14").Value = '[''evolutionary algorithm'', ''fault injection'', ''fitness functions'', ''genetic algorithm'', ''meta-heuristics'', ''mutation testing'', ''programming languages'', ''sbst'', ''search based software testing'', ''simulated annealing'']'
$ws.Range("Q14").Value = 0.07407407407407407
$ws.Range("N15").Value = ""
$ws.Range("O15").Value = '[''anti-patternandexperiment'', ''controlled'', ''disharmony'', ''empirical'', ''ethnography'', ''smell'', ''study'', ''survey'', ''“action research”'', ''“code anomaly”'', ''“design anomaly”'', ''“design flaw”'', ''“exploratory analysis”'']'
$ws.Range("P15").Value = '[''nan'']'
$ws.Range("Q15").Value = 0
$ws.Range("N16").Value = '("Strategic Information Systems" OR "Strategic IS" OR SIA) AND (Inputs OR Data OR Resources OR "Input Factors") AND (Processes OR Methods OR Workflows OR Procedures) AND (Outputs OR Results OR Outcomes OR Deliverables) AND (Usability OR Effectiveness OR "User Experience" OR Adoption)'
$ws.Range("O16").Value = '[''existing systems'', ''legacy'', ''migration'', ''modernization'', ''object-oriented'', ''re-engineering'', ''service identification'', ''service mining'', ''service packaging'', ''transformation'']'
$ws.Range("P16").Value = '[''adoption'', ''data'', ''deliverables'', ''effectiveness'', ''input factors'', ''inputs'', ''methods'', ''outcomes'', ''outputs'', ''procedures'', ''processes'', ''resources'', ''results'', ''sia'', ''strategic information systems'', ''strategic is'', ''usability'', ''user experience'', ''workflows'']'
$ws.Range("Q16").Value = 0
$ws.Range("N17").Value = '("agile method" OR "agile methods" OR "Scrum" OR "XP" OR "Extreme Programming" OR "Kanban") AND ("method tailoring" OR "method adaptation" OR "method customization" OR "method configuration") AND ("research methods" OR "study design" OR "data collection" OR "case studies" OR "implementation" OR "evaluation" OR "context factors" OR "project characteristics" OR "organizational needs")'
$ws.Range("O17").Value = '[''adoption orfdd'', ''adoption orkanban'', ''adoption orlean'', ''adoption orscrum'', ''adoption orxp'', ''adoption or“feature driven development”'', ''agile'', ''method'', ''practice'', ''practice'', ''practice'', ''practice'', ''practice'', ''practice'', ''practice'', ''practice'', ''practice'', ''practice'', ''practice'', ''practice'', ''practice'', ''select'', ''select oragile'', ''select orfdd'', ''select orkanban'', ''select orlean'', ''select orxp'', ''select or“feature driven development”'', ''tailoring orfdd'', ''tailoring orkanban'', ''tailoring orlean'', ''tailoring orscrum'', ''tailoring orscrum'', ''tailoring orxp'', ''tailoring or“feature driven development”'']'
$ws.Range("P17").Value = '[''agile method'', ''agile methods'', ''case studies'', ''context factors'', ''data collection'', ''evaluation'', ''extreme programming'', ''implementation'', ''kanban'', ''method adaptation'', ''method configuration'', ''method customization'', ''method tailoring'', ''organizational needs'', ''project characteristics'', ''research methods'', ''scrum'', ''study design'', ''xp'']'
$ws.Range("Q17").Value = 0
$ws.Range("N18").Value = '("agile practice*" OR "agile method*" OR scrum OR kanban) AND ("critical success factor*" OR CSF) AND ("global software development" OR GSD OR "distributed software development") AND ("green software" OR "sustainable software" OR "green IT" OR "sustainable IT")'
$ws.Range("O18").Value = '[''“agile methods”'', ''“agile”'', ''“green agile”'', ''“green software engineering”'', ''“green software”'', ''“greener software”'', ''“practices”'', ''“solutions”'', ''“sustainable software”'']'
$ws.Range("P18").Value = '[''agile method*'', ''agile practice*'', ''critical success factor*'', ''csf'', ''distributed software development'', ''global software development'', ''green it'', ''green software'', ''gsd'', ''kanban'', ''scrum'', ''sustainable it'', ''sustainable software'']'
$ws.Range("Q18").Value = 0
$ws.Range("N19").Value = '("CMMI-DEV" OR "Capability Maturity Model Integration") AND (Agile OR Scrum OR Kanban OR "Extreme Programming") AND (Web OR "Web-based" OR "Web environment" OR "Web application") AND ("maturity level*" OR "process area*") AND ("evaluation criteria" OR "validation" OR "experiment*" OR "case study")'
$ws.Range("O19").Value = '[''agile'', ''agility'', ''capability maturity model'', ''cmmi'', ''extreme programming'', ''scrum'', ''web'', ''web engineering'']'
$ws.Range("P19").Value = '[''agile'', ''capability maturity model integration'', ''case study'', ''cmmi-dev'', ''evaluation criteria'', ''experiment*'', ''extreme programming'', ''kanban'', ''maturity level*'', ''process area*'', ''scrum'', ''validation'', ''web'', ''web application'', ''web environment'', ''web-based'']'
$ws.Range("Q19").Value = 0.2
$ws.Range("N20").Value = ""
$ws.Range("O20").Value = '[''desired skills'', ''educational needs'', ''essential competencies'', ''knowledge needs'', ''knowledge requirements'', ''skill requirements'', ''software developers'', ''software engineers'']'
$ws.Range("P20").Value = '[''nan'']'
$ws.Range("Q20").Value = 0
$ws.Range("N21").Value = '("app review analysis" OR "application review analysis" OR "user feedback analysis" OR "sentiment analysis" OR "text mining") AND (techniques OR methods OR approaches OR algorithms) AND (software engineering OR software development OR requirements engineering OR maintenance OR testing) AND (empirical evaluation OR validation OR case study OR experiment) AND (effectiveness OR performance OR support OR impact)'
$ws.Range("O21").Value = '["''app review''", "''app store review''", "''app store''", "''requirement engineering''", "''software coding''", "''software configuration''", "''software construction''", "''software design''", "''software development''", "''software engineering''", "''software maintenance''", "''software quality''", "''software requirement''", "''software testing''", "''user feedback''", "''user review''"]'
$ws.Range("P21").Value = '[''algorithms'', ''app review analysis'', ''application review analysis'', ''approaches'', ''case study'', ''effectiveness'', ''empirical evaluation'', ''experiment'', ''impact'', ''maintenance'', ''methods'', ''performance'', ''requirements engineering'', ''sentiment analysis'', ''software development'', ''software engineering'', ''support'', ''techniques'', ''testing'', ''text mining'', ''user feedback analysis'', ''validation'']'
$ws.Range("Q21").Value = 0
$ws.Range("N22").Value = '("offshore development" OR "offshore outsourcing") AND (vendor OR "service provider") AND (reliability OR "service quality" OR dependability) AND (challenges OR issues OR risks OR problems) AND (impact OR effect) AND (client OR customer) AND (prioritization OR "risk assessment" OR "value allocation" OR categorization OR "critical issues")'
$ws.Range("O22").Value = '[''dealer'', ''developer'', ''marketer'', ''merchant'', ''retailer'', ''salesperson'', ''seller'', ''service-provider'', ''trader'', ''wholesaler”'', ''“barriers”'', ''“challenges”'', ''“critical factors”'', ''“global software development”'', ''“hurdles”'', ''“issues”'', ''“obstacles”'', ''“risk analysis”'', ''“risks”'', ''“software outsourcing development”'', ''“software outsourcing”'', ''“vendor'']'
$ws.Range("P22").Value = '[''categorization'', ''challenges'', ''client'', ''critical issues'', ''customer'', ''dependability'', ''effect'', ''impact'', ''issues'', ''offshore development'', ''offshore outsourcing'', ''prioritization'', ''problems'', ''reliability'', ''risk assessment'', ''risks'', ''service provider'', ''service quality'', ''value allocation'', ''vendor'']'
$ws.Range("Q22").Value = 0
$ws.Range("N23").Value = '("software testing" OR "testing domain") AND (ontology OR ontologies) AND (concept OR concepts OR relationship OR relationships OR property OR properties OR axiom OR axioms) AND (classification OR classified OR categorization OR categorized)'
$ws.Range("O23").Value = '[''ontologies'', ''ontology'', ''software test'', ''software testing'']'
$ws.Range("P23").Value = '[''axiom'', ''axioms'', ''categorization'', ''categorized'', ''classification'', ''classified'', ''concept'', ''concepts'', ''ontologies'', ''ontology'', ''properties'', ''property'', ''relationship'', ''relationships'', ''software testing'', ''testing domain'']'
$ws.Range("Q23").Value = 0.1764705882352941
$ws.Range("N24").Value = ""
$ws.Range("O24").Value = '[''“experience curve”'', ''“learning curve”'', ''“software development”'', ''“software engineering”'', ''“software process”'', ''“software project”'']'
$ws.Range("P24").Value = '[''nan'']'
$ws.Range("Q24").Value = 0
$ws.Range("N25").Value = '("IT personnel" OR "information technology staff" OR "computer professionals" OR "IT employees") AND ("intentions to leave" OR "turnover intention" OR "employee attrition" OR "job abandonment" OR "resignation") AND (drivers OR factors OR reasons OR motivations OR causes OR predictors OR determinants) AND (workplace OR organization OR company OR employment)'
$ws.Range("O25").Value = '[''“employee retention”'', ''“employees retention”'', ''“intention for turnover”'', ''“intention to leave”'', ''“intention to quit”'', ''“intention to stay”'', ''“intention to withdraw”'', ''“is employee”'', ''“is manager”'', ''“is personnel”'', ''“is professional”'', ''“is workforce”'', ''“it employee”'', ''“it manager”'', ''“it personnel”'', ''“it professional”'', ''“it worker”'', ''“it workforce”'', ''“leave intention”'', ''“manager retention”'', ''“managers retention”'', ''“personnel retention”'', ''“professional retention”'', ''“professionals retention”and “is worker”'', ''“quit intention”'', ''“software designer”'', ''“software developer”'', ''“software engineer”'', ''“software manager”'', ''“software programmer”'', ''“software project manager”'', ''“stay intention”'', ''“system analyst”'', ''“turnover intention”'', ''“withdrawal intention”'', ''“worker retention”'', ''“workers retention”'']'
$ws.Range("P25").Value = '[''causes'', ''company'', ''computer professionals'', ''determinants'', ''drivers'', ''employee attrition'', ''employment'', ''factors'', ''information technology staff'', ''intentions to leave'', ''it employees'', ''it personnel'', ''job abandonment'', ''motivations'', ''organization'', ''predictors'', ''reasons'', ''resignation'', ''turnover intention'', ''workplace'']'
$ws.Range("Q25").Value = 0
$ws.Range("N26").Value = ""
$ws.Range("O26").Value = '[''agile'', ''hcd'', ''hci'', ''hmi'', ''kanban'', ''lean'', ''scrum'', ''ucd'', ''usability'', ''ux'', ''“design thinking”and“user experience”'', ''“extreme programming”'']'
$ws.Range("P26").Value = '[''nan'']'
$ws.Range("Q26").Value = 0

Write-Output "done"